# Update "想去人数" (want-to-go count) values on sheets "展览" and "全部类型"
# per the published-site regeneration (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 1791
$wsExpo.Range("F6").Value  = 1134
$wsExpo.Range("F7").Value  = 2250
$wsExpo.Range("F8").Value  = 2170
$wsExpo.Range("F12").Value = 1700
$wsExpo.Range("F17").Value = 250
$wsExpo.Range("F18").Value = 1610
$wsExpo.Range("F20").Value = 1329
$wsExpo.Range("F21").Value = 750
$wsExpo.Range("F24").Value = 12398
$wsExpo.Range("F32").Value = 1937

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F10").Value = 53

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 1791
$wsAll.Range("F7").Value  = 1134
$wsAll.Range("F8").Value  = 2250
$wsAll.Range("F9").Value  = 2170
$wsAll.Range("F14").Value = 1700
$wsAll.Range("F22").Value = 250
$wsAll.Range("F23").Value = 1610
$wsAll.Range("F25").Value = 1329
$wsAll.Range("F26").Value = 750
$wsAll.Range("F30").Value = 12398
$wsAll.Range("F40").Value = 1937
$wsAll.Range("F47").Value = 53
